$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# =====================================================================
# 1. Add 20 new rows (182-201) for the new 'Cardiology' domain, which
#    belongs to the 'Medicine & Health' domain cluster (topic #6).
# =====================================================================
$domainCluster = 'Medicine & Health'
$topicNumber = 6
$domainName = 'Cardiology'
$promptLabel = 'Domain_FSPrompt'

$r = 182
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'adversarial learning'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to apply adversarial learning techniques to detect and mitigate potential adversarial attacks on ECG data classifiers, ensuring the reliability and robustness of heart abnormality diagnoses.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 183
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'cnn'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to explore CNN-based techniques for real-time analysis of wearable device data, such as continuous heart rate monitoring, to detect anomalies and provide timely alerts for patients with underlying cardiac conditions.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 184
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'conversational agent'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to integrate a conversational agent into electronic health records (EHR) systems, enabling patients to easily access and update their medical histories and symptoms, facilitating more efficient and accurate clinical assessments.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 185
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'decision tree'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to construct decision tree models using patient demographics and medical history to predict the likelihood of adverse cardiac events within the next five years, aiding in early intervention and risk management strategies.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 186
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'document classification'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to develop a document classification system to classify electronic health records (EHRs) based on patient symptoms and diagnostic tests, aiding in the identification of patterns and trends in cardiovascular diseases.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 187
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'entity extraction'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to develop an entity extraction system to automatically identify and extract key cardiac parameters (e.g., ejection fraction, QT interval) from clinical notes and reports, enabling faster analysis and decision-making.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 188
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'feature selection'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to apply feature selection methods to filter out irrelevant or redundant features from ECG data, so that I can enhance the performance of algorithms detecting cardiac arrhythmias.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 189
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'imbalanced dataset'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to use keyword extraction algorithms to extract relevant terms from medical literature and clinical guidelines pertaining to cardiac rehabilitation protocols, aiding in the development of evidence-based treatment plans.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 190
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'keyword extraction'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to address class imbalance in datasets used for predicting rare cardiac conditions using machine learning algorithms, ensuring accurate identification and early intervention for patients at higher risk.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 191
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'k-nearest neighbor'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to use k-NN models to predict patient-specific responses to different cardiac medications based on similar patient profiles, facilitating personalized treatment plans for heart disease management.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 192
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'multi-label classification'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to implement multi-label classification algorithms to classify cardiac imaging studies (e.g., echocardiograms, CT scans) into multiple diagnostic categories (e.g., valve disease, myocardial infarction), aiding in accurate disease characterization.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 193
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'neural network'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to develop a neural network-based system to predict the progression of heart failure in patients based on dynamic changes in biomarkers and clinical indicators, guiding timely interventions and patient monitoring.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 194
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'random forest'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to deploy a random forest model to predict the likelihood of adverse drug reactions in cardiac patients based on medication history, comorbidities, and genetic predispositions, optimizing medication management strategies.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 195
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'semantic similarity'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to develop a semantic similarity model to compare clinical notes and identify similar cases of coronary artery disease, aiding in pattern recognition and treatment planning.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 196
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'sentiment analysis'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to implement sentiment analysis on patient feedback from cardiac rehabilitation programs to assess overall patient satisfaction and identify areas for program improvement.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 197
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'speech to text'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to develop a speech to text system to transcribe cardiology consultations and patient histories accurately, improving documentation efficiency and clinical workflow.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 198
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'text categorization'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to develop a text categorization model to classify medical literature and research articles on various cardiac conditions (e.g., myocardial infarction, arrhythmias) for easier access and knowledge synthesis.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 199
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'unsupervised clustering'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to apply unsupervised clustering algorithms to group patients based on similar risk factor profiles (e.g., smoking history, cholesterol levels), enabling targeted preventive interventions for cardiovascular diseases.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 200
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'voice recognition'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to develop a voice recognition system to accurately transcribe cardiology consultations and patient histories from audio recordings, improving documentation accuracy and efficiency.'
$ws.Cells.Item($r, 6).Value = $promptLabel

$r = 201
$ws.Cells.Item($r, 1).Value = $domainCluster
$ws.Cells.Item($r, 2).Value = $topicNumber
$ws.Cells.Item($r, 3).Value = $domainName
$ws.Cells.Item($r, 4).Value = 'word embedding'
$ws.Cells.Item($r, 5).Value = 'As a cardiologist, I want to utilize word embedding techniques to represent clinical terms and medical concepts from cardiology literature, enabling more accurate semantic understanding and retrieval of relevant research findings.'
$ws.Cells.Item($r, 6).Value = $promptLabel

# =====================================================================
# 2. Highlight the Domain cluster/Topic/Domain columns (A:C) of the new
#    rows with the same light-blue background used for this new block.
#    A clean helper cell is formatted first and its format is then
#    copied onto the target range so that only one new fill/style is
#    created (mirrors how the other domain blocks were color-coded).
# =====================================================================
$helper = $ws.Range("H205")
$helper.Interior.Color = 15441517
$newRowsRange = $ws.Range("A182:C201")
$helper.Copy()
$newRowsRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$helper.Clear()

# =====================================================================
# 3. Normalize the formatting of rows 2-41 (the first domain block) so
#    the User Story column (E) carries the same explicit style as the
#    rest of the sheet, and let Excel recompute the (no-longer custom)
#    row heights.
# =====================================================================
$formatSource = $ws.Range("E42")
$formatTarget = $ws.Range("E2:E41")
$formatSource.Copy()
$formatTarget.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows("2:41").AutoFit()

# =====================================================================
# 4. Update the sheet view so it matches the position left after the
#    edit (scrolled down to the newly added rows).
# =====================================================================
$ws.Range("D199").Select()
$excel.ActiveWindow.ScrollRow = 175

